$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Testmail #1: Kun jij dit even regelen?"
$ws.Range("B6").Value = "Beste klant,`nBedankt voor je e-mail. Om je vraag beter te kunnen beantwoorden, heb ik meer details nodig. Kunt u beschrijven waar u specifiek hulp bij nodig heeft? Als u meer informatie geeft, kan ik u beter van dienst zijn.`nMet vriendelijke groet,`n[Naam]`nE-mailassistent"
$ws.Range("C6").Value = "Kun jij dit even regelen?"
$ws.Range("D6").Value = "mailmind.test@zohomail.eu"
$ws.Range("E6").Value = "Overig"
$ws.Range("F6").Value = "2025-08-01 23:16:44"
$ws.Range("G6").Value = "Ja"
$ws.Range("H6").Value = "Nee"
$ws.Range("I6").Value = "Ja"
$ws.Range("J6").Value = "Nee"

$ws.Rows.Item(6).EntireRow.AutoFit()
